# Update LAVRAS_DO_SUL.xlsx:
#  - Rename "Paineis DARQ" -> "PAINEIS DARQ"
#  - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#  - Delete the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wb.Worksheets("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Keep the first sheet ("PAINEIS DARQ") the active/selected tab, same as before the edit.
$wb.Worksheets("PAINEIS DARQ").Activate()
